$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old vessel data rows (rows 2-10), columns A and B
$ws.Range("A2:B10").ClearContents()

# New header row values across A1:Q1
$headers = @(
    "PORT",
    "MICT SERVICE NAME",
    "SERVICE NAME",
    "SERVICE DESC",
    "ROUTE",
    "LEAD SL",
    "SAILING FREQ",
    "PARTICIPANTS",
    "VESSEL OPERATOR",
    "# OF VESSELS",
    "# OF VESSELS PER ROW COUNT",
    "WEEKLY CAPACITY",
    "SHIPS USED",
    "PORT ROTATION",
    "ALT SRVC CD",
    "VESSEL SIZE",
    "VESSEL_NAME"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
